$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "To" time for the last entry (row 14) from 20:00 to 22:00.
# This is a time value stored as a fraction of a day: 22/24 = 0.916666666666667
$ws.Range("C14").Value = 22 / 24

# Move the active selection to C15, matching the saved selection state.
$ws.Range("C15").Select()
